# Weekly update for "Hortaliza, Vega Modelo de Temuco - Acelga":
# a new observation is inserted at the top of the data table (row 324),
# pushing the existing rows 324-353 down to 325-354.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 324; everything from 324 down shifts to 325+.
$ws.Rows.Item(324).Insert()

# Populate the new row 324 with the latest weekly record.
$ws.Cells.Item(324, 1).Value  = 10
$ws.Cells.Item(324, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(324, 3).Value  = "La Araucanía"
$ws.Cells.Item(324, 4).Value  = 44769
$ws.Cells.Item(324, 5).Value  = 9
$ws.Cells.Item(324, 6).Value  = 100112009
$ws.Cells.Item(324, 7).Value  = "Acelga"
$ws.Cells.Item(324, 8).Value  = "Sin especificar"
$ws.Cells.Item(324, 9).Value  = "Primera"
$ws.Cells.Item(324, 10).Value = 200
$ws.Cells.Item(324, 11).Value = 9000
$ws.Cells.Item(324, 12).Value = 9000
$ws.Cells.Item(324, 13).Value = 9000
$ws.Cells.Item(324, 14).Value = '$/docena de atados (12 kilos)'
$ws.Cells.Item(324, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(324, 16).Value = 750
$ws.Cells.Item(324, 17).Value = 12
$ws.Cells.Item(324, 18).Value = "Hortaliza"
